# Chiffres COVID-19 Valais - daily data update
# Revises a few previously-reported new-case counts and fills in the
# previously-empty row for 2021-03-02 (row 371) now that data is available.
# Columns B, H, J and K are shared formulas and recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: columns L and M are formatted as Text (@), so a plain .Value
# assignment would be stored as a text string instead of a number (this
# matches real Excel's "typing a number into a Text cell" behaviour).
# Temporarily switch the cell to General, write the number, then restore
# the original (Text) number format so the cell keeps its look/format.
function Set-NumericValue($range, $value) {
    $fmt = $range.NumberFormat
    $range.NumberFormat = "General"
    $range.Value = $value
    $range.NumberFormat = $fmt
}

# Row 367 (2021-02-26): new-case count revised down from 69 to 68
$ws.Range("C367").Value = 68

# Row 369 (2021-02-28): new-case count revised from 25 to 34
$ws.Range("C369").Value = 34

# Row 370 (2021-03-01): new-case count revised from 14 to 64
$ws.Range("C370").Value = 64

# Row 371 (2021-03-02): previously blank (no data yet) - now filled in
$ws.Range("C371").Value = 6
$ws.Range("E371").Value = 8
$ws.Range("F371").Value = 7
$ws.Range("G371").Value = 25
Set-NumericValue $ws.Range("L371") 0
Set-NumericValue $ws.Range("M371") 0
